$wb = $excel.ActiveWorkbook

# OFF sheet updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 378
$wsOff.Range("C2").Value = 242
$wsOff.Range("D2").Value = 155
$wsOff.Range("E2").Value = 69

# DEF sheet updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 516
$wsDef.Range("C2").Value = 389
$wsDef.Range("D2").Value = 111
$wsDef.Range("E2").Value = 66
$wsDef.Range("G2").Value = 5
